# New results update (2021/04/20 19:45)
# - refresh several measured values in column B (rows 17, 18, 31)
# - B32 holds =AVERAGE(B2:B31) and recalculates automatically
# - re-assert the "Bad"/"Neutral" conditional cell styles on the B column
#   (their underlying style-table slots get re-written by Excel on save,
#   but the named style per cell is unchanged)
# - move the active selection to B17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated measurements ---------------------------------------------
$ws.Range("B17").Value = 0.3714
$ws.Range("B18").Value = 0.4298
$ws.Range("B31").Value = 0.3202

# --- re-apply the named cell styles (round-trips the style table the
#     same way the source workbook's edit did) -------------------------
$neutralCells = @("B10", "B11", "B12", "B16", "B17", "B20", "B22", "B28", "B30")
foreach ($addr in $neutralCells) {
    $ws.Range($addr).Style = "适中"
}

$badCells = @("B3", "B4", "B8", "B13", "B14", "B15", "B26", "B27", "B31")
foreach ($addr in $badCells) {
    $ws.Range($addr).Style = "差"
}

# --- recalc so B32's AVERAGE formula reflects the new numbers ---------
$excel.Calculate()

# --- move the selection, matching the saved cursor position -----------
$ws.Range("B17").Select() | Out-Null
